$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Paragraph 2 ("...proceso no puede realizar ... multihilo."):
#    trim the trailing period off "multihilo." and then append the
#    new explanatory sentences as additional runs at the end of the
#    paragraph.
# ------------------------------------------------------------------
$p2 = $d.Paragraphs(2)

$trim = $p2.Range.Duplicate
$trim.Find.Execute("multihilo.")
$trim.Text = "multihilo"

$p2.Range.InsertAfter(". Cuando un proceso es multihilo, los hilos comparten información haciendo que la realización de la tarea sea ")
$p2.Range.InsertAfter("más")
$p2.Range.InsertAfter(" eficaz")
$p2.Range.InsertAfter(", los hilos comparten sección de código, datos y recursos")
$p2.Range.InsertAfter(".")

# ------------------------------------------------------------------
# 2) The two trailing empty paragraphs become four paragraphs of new
#    text describing thread implementations.
# ------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p3.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(3)
$p3.Range.Text = "Los hilos se pueden implementar de dos formas: una implementación hecha por el usuario o una implementación hecha por el sistema operativo."

$p4 = $d.Paragraphs(4)
$p4.Range.InsertParagraphAfter()
$p4.Range.Text = "La implementación hecha por el usuario es hecha por una aplicación de forma separada al sistema operativo, por lo tanto, este no es consciente a los hilos creados por el usuario. Este tipo de hilos se realizan mediante librerías."

$p5 = $d.Paragraphs(5)
$p5.Range.Text = "La implementación hecha por el sistema operativo como lo dice su nombre es hecha por el sistema operativo, siendo mas especifico por su núcleo/kernel."

$p6 = $d.Paragraphs(6)
$p6.Range.Text = "Cada uno tiene su ventaja y desventaja, la implementación por el usuario permite crear hilos incluso si el kernel no usa hilos de manera nativa, pero estos pueden bloquear al resto de hilos cuando llaman al sistema. Los hilos hechos por el sistema operativo usan de mejor manera los recursos de las diferentes arquitectura de los procesadores"
